# Regional Availability Factor.xlsx -- "updated 4.0 files and mdl"
#
# 1) About!C1 source date bumped two weeks forward (45366 -> 45379, i.e. 3/15/2024 -> 3/28/2024)
# 2) RAF-capacity!B24 and B25 (hydrogen combustion turbine / hydrogen combined cycle
#    capacity-credit multipliers) raised from 0.3 to 1
# 3) The workbook was left with RAF-capacity as the active/selected sheet, scrolled down
#    and zoomed to 80%, with B25 selected.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "last updated" source date on the About sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = "3/28/2024"

# --- 2. Update the RAF-capacity hydrogen plant capacity-credit multipliers ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Column A on RAF-capacity picked up an explicit (slightly narrower) width
$wsCapacity.Columns.Item(1).ColumnWidth = 28.1

# --- 3. Leave the workbook focused on RAF-capacity, matching the saved view state ---
[void]$wsCapacity.Select()
$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 14
$win.ScrollColumn = 1
[void]$wsCapacity.Range("B25").Select()
